$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 24 reconciliation rows (rows 163-186) with data for columns A:I
$newRows = New-Object 'object[,]' 24,9

$newRows[0,0] = 237681102046
$newRows[0,1] = "FRANCOISE NKENFACK NKENGMO"
$newRows[0,2] = "Rte_3"
$newRows[0,3] = "Essec"
$newRows[0,4] = 84169.8
$newRows[0,5] = 397988
$newRows[0,6] = 313818.2
$newRows[0,7] = 4.728394269678673
$newRows[0,8] = "Cite Sic"

$newRows[1,0] = 237682368679
$newRows[1,1] = "MFS SIM AA 2"
$newRows[1,2] = "Rte_2"
$newRows[1,3] = "Essec"
$newRows[1,4] = 219379.0363636363
$newRows[1,5] = 322045
$newRows[1,6] = 102665.9636363637
$newRows[1,7] = 1.467984386011194
$newRows[1,8] = "Cite Sic"

$newRows[2,0] = 237683360459
$newRows[2,1] = "LUCIE MAJOLIE LELE NKANKEU"
$newRows[2,2] = "Rte_0"
$newRows[2,3] = "Essec"
$newRows[2,4] = 5000
$newRows[2,5] = 574
$newRows[2,6] = -4426
$newRows[2,7] = 0.1148
$newRows[2,8] = "Cite Sic"

$newRows[3,0] = 237652899422
$newRows[3,1] = "NOUMOU epouse SAGNON MARCELINE LA NEGRESSE"
$newRows[3,2] = "Rte_0"
$newRows[3,3] = "Hopital General Douala"
$newRows[3,4] = 68881.36363636363
$newRows[3,5] = 155671
$newRows[3,6] = 86789.63636363637
$newRows[3,7] = 2.259987198015033
$newRows[3,8] = "Ndogbong"

$newRows[4,0] = 237670904526
$newRows[4,1] = "MFS SIM PROVISOIRE 20"
$newRows[4,2] = "Rte_5"
$newRows[4,3] = "Hopital General Douala"
$newRows[4,4] = 352072.5
$newRows[4,5] = 0
$newRows[4,6] = -352072.5
$newRows[4,7] = 0
$newRows[4,8] = "Ndogbong"

$newRows[5,0] = 237671105116
$newRows[5,1] = "MFS  AM FACE HÔPITAL GÉNÉRAL"
$newRows[5,2] = "Rte_3"
$newRows[5,3] = "Hopital General Douala"
$newRows[5,4] = 112588.4615384615
$newRows[5,5] = 18
$newRows[5,6] = -112570.4615384615
$newRows[5,7] = 0.0001598742868855259
$newRows[5,8] = "Ndogbong"

$newRows[6,0] = 237672916354
$newRows[6,1] = "MAMADOU DIAN BAH LENA GLOBAL"
$newRows[6,2] = "Rte_0"
$newRows[6,3] = "Hopital General Douala"
$newRows[6,4] = 24826.15384615385
$newRows[6,5] = 68564
$newRows[6,6] = 43737.84615384616
$newRows[6,7] = 2.761764888145256
$newRows[6,8] = "Ndogbong"

$newRows[7,0] = 237672920086
$newRows[7,1] = "NAMY NGOKO CLARISSE ROSE VERTINE KAMILAH CONNECTION"
$newRows[7,2] = "Rte_5"
$newRows[7,3] = "Hopital General Douala"
$newRows[7,4] = 113740
$newRows[7,5] = 6730
$newRows[7,6] = -107010
$newRows[7,7] = 0.05917003692632319
$newRows[7,8] = "Ndogbong"

$newRows[8,0] = 237674000053
$newRows[8,1] = "FRED JUNIOR ZOK EDOU"
$newRows[8,2] = "Rte_5"
$newRows[8,3] = "Hopital General Douala"
$newRows[8,4] = 56545
$newRows[8,5] = 138362
$newRows[8,6] = 81817
$newRows[8,7] = 2.446936068617915
$newRows[8,8] = "Ndogbong"

$newRows[9,0] = 237674841555
$newRows[9,1] = "BEATRICE TCHAMTIEU EPSE NGAMENI"
$newRows[9,2] = "Rte_5"
$newRows[9,3] = "Hopital General Douala"
$newRows[9,4] = 150379.6
$newRows[9,5] = 192473
$newRows[9,6] = 42093.39999999999
$newRows[9,7] = 1.279914296886014
$newRows[9,8] = "Ndogbong"

$newRows[10,0] = 237674899678
$newRows[10,1] = "VIVIANE MADJUIMEKEM FOMEKONG"
$newRows[10,2] = "Rte_5"
$newRows[10,3] = "Hopital General Douala"
$newRows[10,4] = 159035.175
$newRows[10,5] = 82338
$newRows[10,6] = -76697.17499999999
$newRows[10,7] = 0.517734520051932
$newRows[10,8] = "Ndogbong"

$newRows[11,0] = 237676439452
$newRows[11,1] = "SAGNOU BRINDA JOSELINE _DIGITAL BUSINESS SARL"
$newRows[11,2] = "Rte_5"
$newRows[11,3] = "Hopital General Douala"
$newRows[11,4] = 87377.35
$newRows[11,5] = 229
$newRows[11,6] = -87148.35
$newRows[11,7] = 0.002620816493061417
$newRows[11,8] = "Ndogbong"

$newRows[12,0] = 237676695935
$newRows[12,1] = "SIMON PIERRE AKOA"
$newRows[12,2] = "Rte_5"
$newRows[12,3] = "Hopital General Douala"
$newRows[12,4] = 88445
$newRows[12,5] = 132533
$newRows[12,6] = 44088
$newRows[12,7] = 1.498479280909039
$newRows[12,8] = "Ndogbong"

$newRows[13,0] = 237677745809
$newRows[13,1] = "FRANKLIN MUA ZUO"
$newRows[13,2] = "Rte_0"
$newRows[13,3] = "Hopital General Douala"
$newRows[13,4] = 23439.6
$newRows[13,5] = 228443
$newRows[13,6] = 205003.4
$newRows[13,7] = 9.74602808921654
$newRows[13,8] = "Ndogbong"

$newRows[14,0] = 237679127464
$newRows[14,1] = "DJUFFO TSOATA MARIE NOEL KAMILAH CONNECTION GROUP"
$newRows[14,2] = "Rte_0"
$newRows[14,3] = "Hopital General Douala"
$newRows[14,4] = 16950
$newRows[14,5] = 33636
$newRows[14,6] = 16686
$newRows[14,7] = 1.984424778761062
$newRows[14,8] = "Ndogbong"

$newRows[15,0] = 237679422291
$newRows[15,1] = "ETS LE CONTENT 32"
$newRows[15,2] = "Rte_0"
$newRows[15,3] = "Hopital General Douala"
$newRows[15,4] = 101597.1428571429
$newRows[15,5] = 100022
$newRows[15,6] = -1575.142857142855
$newRows[15,7] = 0.9844961894316488
$newRows[15,8] = "Ndogbong"

$newRows[16,0] = 237651433330
$newRows[16,1] = "NGUIAZONG DORIANE LAURE KAMILAH CONNECTION GROUP"
$newRows[16,2] = "Rte_7"
$newRows[16,3] = "Makepe Conquete"
$newRows[16,4] = 59904.58333333334
$newRows[16,5] = 83422
$newRows[16,6] = 23517.41666666666
$newRows[16,7] = 1.392581257694528
$newRows[16,8] = "Ndogbong"

$newRows[17,0] = 237654168696
$newRows[17,1] = "DZEUMAZONG FLORENCE ETS MOBILE FINANCIAL SERVICES MFS"
$newRows[17,2] = "Rte_0"
$newRows[17,3] = "Makepe Conquete"
$newRows[17,4] = 7797.028571428572
$newRows[17,5] = 1893
$newRows[17,6] = -5904.028571428572
$newRows[17,7] = 0.2427847971388367
$newRows[17,8] = "Ndogbong"

$newRows[18,0] = 237670799877
$newRows[18,1] = "ETS MOBILE FINANCIAL SERVICES MFS MENANDJIO HORTENSE BIENVENUE"
$newRows[18,2] = "Rte_0"
$newRows[18,3] = "Makepe Conquete"
$newRows[18,4] = 23610
$newRows[18,5] = 642657
$newRows[18,6] = 619047
$newRows[18,7] = 27.21969504447268
$newRows[18,8] = "Ndogbong"

$newRows[19,0] = 237671351291
$newRows[19,1] = "MFS LTDLA CBOX R3 MOUTHIEU JOSETTE CHANCELINE"
$newRows[19,2] = "Rte_5"
$newRows[19,3] = "Makepe Conquete"
$newRows[19,4] = 95869.86000000002
$newRows[19,5] = 324495
$newRows[19,6] = 228625.14
$newRows[19,7] = 3.384744694526517
$newRows[19,8] = "Ndogbong"

$newRows[20,0] = 237671378136
$newRows[20,1] = "KOUBINOM DIPITA SARIETTE CRISTELLE ETS MOBILE FINANCIAL SERVICES MFS"
$newRows[20,2] = "Rte_6"
$newRows[20,3] = "Makepe Conquete"
$newRows[20,4] = 54416.5
$newRows[20,5] = 4426
$newRows[20,6] = -49990.5
$newRows[20,7] = 0.08133562430512804
$newRows[20,8] = "Ndogbong"

$newRows[21,0] = 237671605749
$newRows[21,1] = "ETS TCHATCHOUANG PAUL  ETP LTDLA CBOX RO MEGAPTCHE VICTORINE"
$newRows[21,2] = "Rte_6"
$newRows[21,3] = "Makepe Conquete"
$newRows[21,4] = 107695
$newRows[21,5] = 436933
$newRows[21,6] = 329238
$newRows[21,7] = 4.057133571660708
$newRows[21,8] = "Ndogbong"

$newRows[22,0] = 237671615641
$newRows[22,1] = "BEGO FOGUE CHRISTELLE KAMILAH CONNECTION GROUP"
$newRows[22,2] = "Rte_0"
$newRows[22,3] = "Makepe Conquete"
$newRows[22,4] = 8701.666666666666
$newRows[22,5] = 21212
$newRows[22,6] = 12510.33333333333
$newRows[22,7] = 2.437693928366214
$newRows[22,8] = "Ndogbong"

$newRows[23,0] = 237673739931
$newRows[23,1] = "MOFFO GERMAIN SPECTRUM SPECTRUM"
$newRows[23,2] = "Rte_0"
$newRows[23,3] = "Makepe Conquete"
$newRows[23,4] = 21375
$newRows[23,5] = 45501
$newRows[23,6] = 24126
$newRows[23,7] = 2.128701754385965
$newRows[23,8] = "Ndogbong"

$startCell = $ws.Cells.Item(163, 1)
$endCell = $ws.Cells.Item(186, 9)
$targetRange = $ws.Range($startCell, $endCell)
$targetRange.Value = $newRows

Write-Output "Wrote rows $($startCell.Row) to $($endCell.Row)"
